# Rename the existing sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Audi TTS"

# Add a new worksheet right after "Audi TTS" and rename it
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Audi TT Sport"

# Move the "Audi TT Sport" rows (rows 9 and 10) from ws1 into ws2
# Header row
$ws1.Range("A1:J1").Copy()
$ws2.Range("A1").PasteSpecial(-4104)

# Data rows (rows 9-10 on ws1 map to rows 2-3 on ws2)
$ws1.Range("A9:J10").Copy()
$ws2.Range("A2").PasteSpecial(-4104)

$excel.CutCopyMode = 0

# Remove the moved rows from ws1 without shifting subsequent rows up
$ws1.Range("A9:J10").Clear()

# Fix up the selections to match the final state
$ws1.Range("A1:J1").Select()
$ws2.Range("A1:J1").Select()
$ws2.Activate()

Write-Host "Sheet count:" $wb.Worksheets.Count
Write-Host "Sheet1 name:" $ws1.Name
Write-Host "Sheet2 name:" $ws2.Name
